# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6654
$ws1.Range("F7").Value = 543
$ws1.Range("F11").Value = 3
$ws1.Range("F14").Value = 1288
$ws1.Range("F16").Value = 3319
$ws1.Range("F18").Value = 213
$ws1.Range("F19").Value = 1964
$ws1.Range("F20").Value = 76

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6654
$ws4.Range("F8").Value = 543
$ws4.Range("F12").Value = 3
$ws4.Range("F15").Value = 1288
$ws4.Range("F17").Value = 3319
$ws4.Range("F19").Value = 213
$ws4.Range("F20").Value = 1964
$ws4.Range("F21").Value = 76
